$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 174
$ws.Range("F6").Value = 1051
$ws.Range("F7").Value = 1033
$ws.Range("F8").Value = 8006
$ws.Range("F9").Value = 129
$ws.Range("F10").Value = 193
$ws.Range("F11").Value = 6838
$ws.Range("F12").Value = 161
$ws.Range("F14").Value = 4893
$ws.Range("F17").Value = 5310
$ws.Range("F20").Value = 318
$ws.Range("F21").Value = 437
$ws.Range("F22").Value = 308
$ws.Range("F23").Value = 253
$ws.Range("F25").Value = 151
$ws.Range("F26").Value = 95
$ws.Range("F27").Value = 9002
$ws.Range("F29").Value = 1610
$ws.Range("F33").Value = 827
$ws.Range("F35").Value = 72
$ws.Range("F36").Value = 1003
$ws.Range("F37").Value = 1151
$ws.Range("F38").Value = 51
$ws.Range("F39").Value = 4688
$ws.Range("F40").Value = 28
$ws.Range("F42").Value = 1155
$ws.Range("F44").Value = 142
$ws.Range("F45").Value = 71
$ws.Range("F47").Value = 1233
$ws.Range("F48").Value = 31
$ws.Range("F49").Value = 59

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 174
$ws.Range("F8").Value = 1051
$ws.Range("F9").Value = 1033
$ws.Range("F10").Value = 8006
$ws.Range("F11").Value = 129
$ws.Range("F12").Value = 193
$ws.Range("F13").Value = 6838
$ws.Range("F14").Value = 161
$ws.Range("F17").Value = 4893
$ws.Range("F19").Value = 5310
$ws.Range("F21").Value = 316
$ws.Range("F22").Value = 318
$ws.Range("F23").Value = 437
$ws.Range("F24").Value = 308
$ws.Range("F25").Value = 253
$ws.Range("F27").Value = 151
$ws.Range("F28").Value = 95
$ws.Range("F30").Value = 9002
$ws.Range("F32").Value = 1610
$ws.Range("F35").Value = 827
$ws.Range("F37").Value = 72
$ws.Range("F38").Value = 1003
$ws.Range("F39").Value = 1151
$ws.Range("F40").Value = 51
$ws.Range("F41").Value = 4688
$ws.Range("F43").Value = 1155
$ws.Range("F44").Value = 142
$ws.Range("F45").Value = 71
$ws.Range("F47").Value = 1233
$ws.Range("F48").Value = 31
$ws.Range("F49").Value = 59

